$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D:E value writes are stored as text (matching the source data
# which is inline text, not numbers) rather than being auto-converted to
# numeric values by Excel's smart-entry parsing.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.165.67"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.878.13"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "314.01"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "0.5132"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").Value = "0.3912"
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").Value = "0.08356"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "1.120"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("D11").Value = "41.53"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "6.226"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.890.06"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "20.67"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "7.258"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "91.10"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "0.06675"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").Value = "17.79"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "6.045"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "28.215.18"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "11.16"
$ws.Range("D25").Value = "2.267"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").Value = "2.092.91"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "159.73"
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").Value = "20.67"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "125.16"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "0.1060"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "1.037"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "5.852"
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("D34").Value = "3.611"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").Value = "9.669"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").Value = "0.02446"
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("D37").Value = "0.06564"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Value = "0.2188"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").Value = "1.200"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "0.6495"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "4.997"
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("D42").Value = "1.224"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "11.33"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "0.6139"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "13.04"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "1.282"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").Value = "3.676"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "2.021"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").Value = "1.230"
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").Value = "120.76"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "78.20"
$ws.Range("E51").Value = "  -0.66%  "

# Restore the cells to their original (unstyled/default) formatting so
# only the values themselves change.
$dataRange.ClearFormats()
